$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-row unique Name/MiddleName/LastName sample values with a
# single repeated "UserN" placeholder across columns A-D for each row,
# matching the new shared Username value for that row.
$ws.Range("A2:D2").Value = "User1"
$ws.Range("A3:D3").Value = "User2"
$ws.Range("A4:D4").Value = "User3"

# Update the selection to span the whole data block (A2:D4) with A2 active.
$ws.Range("A2:D4").Select()
